$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank rows (rows 6-13) before the old "Run By" row, pushing
# "Run By"/"Run Date" down to rows 14/15 and the header row down from 9 to 17.
# The new blank rows inherit the style of the row above (s=5), matching the
# target formatting for the new label rows.
$ws.Range("A6:A13").EntireRow.Insert()

# Relabel the existing rows that now carry new meaning because new label
# rows were spliced in above the old "Activity Group"/"Activity" rows.
$ws.Range("A4").Value = "Charge Type"
$ws.Range("A5").Value = "Org"

# Fill in the newly inserted label rows.
$ws.Range("A6").Value = "Sector"
$ws.Range("A7").Value = "Subsector"
$ws.Range("A8").Value = "Division"
$ws.Range("A9").Value = "Section"
$ws.Range("A10").Value = "Budget Method"
$ws.Range("A11").Value = "Investment Asset"
$ws.Range("A12").Value = "Activity Group"
$ws.Range("A13").Value = "Activity"

# Reset the active selection to A1 (was C5).
$ws.Range("A1").Select()
